$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (B, C, D, E) keep their text format so numeric-looking
# strings such as "669.37" or "1.00" are not auto-converted to numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.301.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.33%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.797.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.47%  "
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "669.37"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.93%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.97%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.796.17"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.48%  "
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.87%  "
# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.06%  "
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.41%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.55%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.00%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.63"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.99%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.436.85"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.62%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.802.65"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.53%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.270.28"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.42%  "
# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.52%  "
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.13%  "
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.50%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.45"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +19.31%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "475.03"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.07%  "
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.36%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.21"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.17%  "
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.47%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.55%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.29"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.46%  "
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.54%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.947.90"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.62%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.83"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.73%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.31"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.39%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.42"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.41%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.53"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.16%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.178"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +8.45%  "
# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.11"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.14%  "
# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.04%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.753.92"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.71%  "
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.24%  "
# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.51%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.95"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.06%  "
# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.962"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.64%  "
# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.10%  "
# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +10.57%  "
# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.03%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.53"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.29%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.45"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.06%  "
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.74%  "
# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.300"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.48%  "
# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.42"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.83%  "
# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.94%  "
